$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "X data" rows appended for the year cycle 1->12 (indices 204..215),
# extending the normalized feature column used for ML training.
$newData = @(
    @("204", "2.312964634635743E-18"),
    @("205", "3.784851220313033E-17"),
    @("206", "-2.775557561562891E-17"),
    @("207", "-2.775557561562891E-17"),
    @("208", "0"),
    @("209", "0"),
    @("210", "4.163336342344337E-17"),
    @("211", "-2.775557561562892E-18"),
    @("212", "0"),
    @("213", "-3.700743415417188E-17"),
    @("214", "0"),
    @("215", "0")
)

$startRow = 206
$lastFormattedRow = 205

for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $pair = $newData[$i]
    $aVal = [double]$pair[0]
    $bVal = [double]$pair[1]

    # Copy the formatting (bold font, border, alignment) used by the
    # existing index column cells down onto the newly appended row.
    $ws.Range("A$lastFormattedRow").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Range("A$row").Value = $aVal
    $ws.Range("B$row").Value = $bVal
}

$excel.CutCopyMode = 0
